# Update "想去人数" (number of people interested) figures that changed
# between scrapes, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 362
$wsExhibit.Range("F3").Value = 771
$wsExhibit.Range("F4").Value = 268
$wsExhibit.Range("F5").Value = 813
$wsExhibit.Range("F6").Value = 2004
$wsExhibit.Range("F7").Value = 180

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 362
$wsAll.Range("F3").Value = 771
$wsAll.Range("F4").Value = 268
$wsAll.Range("F7").Value = 813
$wsAll.Range("F8").Value = 2004
$wsAll.Range("F10").Value = 180
